$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet
$ws.Name = "Through 2022-10-22"

# Update header text for total column (I1)
$ws.Range("I1").Value = "2022 (through 10-22)"

# Update data values
$ws.Range("I11").Value = 77
$ws.Range("I14").Value = 1354
